$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "1. Click Menu 'Đơn hàng của tôi'`n2. Check danh sách đơn`n3. Click nút 'Chi tiết' đơn đầu tiên`n4. Check thông tin trang chi tiết`n5. Quay lại"
$ws.Range("E2").Value = "Hiển thị đúng danh sách, thông tin chi tiết trùng khớp với tổng tiền, và quay lại thành công"
$ws.Range("F2").Value = "Hoàn thành test: Xem danh sách, chi tiết và quay lại thành công cho đơn #2071"

# Columns resize (bestFit) after the longer text was entered above.
# (Engine quantizes ColumnWidth to 1/6-character steps, so these inputs are
# chosen to land on the stored-width bucket nearest the target.)
$ws.Columns.Item(3).ColumnWidth = 30
$ws.Columns.Item(5).ColumnWidth = 81
$ws.Columns.Item(6).ColumnWidth = 70.83333333333333
